# Missile.xlsx — support more than 1 frame per missile sprite.
# Adds two new table columns (FrameCount, FrameTime) and two new missile
# rows (darkwheel / dragonball), and backfills FrameCount=1 / FrameTime=1
# for the pre-existing single-frame missiles. Also corrects the Image
# (path) offsets for the existing "arrow" rows now that Image values are
# spaced by 10 to make room for multi-frame sprite sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Add the two new table columns (this also sets the row-3 header
#        cells F3/G3 and keeps the table's column <-> header mapping in
#        sync, unlike writing straight into the cells). ------------------
$c6 = $lo.ListColumns.Add()
$c6.Range.Item(1, 1).Value = "FrameCount"
$c7 = $lo.ListColumns.Add()
$c7.Range.Item(1, 1).Value = "FrameTime"

# --- 2. Extend the two label rows above the table (Chinese display name
#        in row 1, field type in row 2) onto the new columns, copying the
#        existing formatting from column E so the header styling matches.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("F2:G2").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

$ws.Range("F1").Value = "图片数"
$ws.Range("G1").Value = "每帧的时间"
$ws.Range("F2").Value = "int"
$ws.Range("G2").Value = "int"

# --- 3. Backfill FrameCount/FrameTime = 1 for the existing single-frame
#        rows, and bump the Image offsets now that each sprite reserves a
#        block of 10. ----------------------------------------------------
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1

$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1

$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1

# --- 4. New missile rows: darkwheel (黑暗轮) and dragonball (黑龙波),
#        each rendered from a 2-frame sprite sheet, 3 time units/frame. --
$ws.Range("A7").Value = 3
$ws.Range("B7").Value = "darkwheel"
$ws.Range("C7").Value = "黑暗轮"
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 30
$ws.Range("F7").Value = 2
$ws.Range("G7").Value = 3

$ws.Range("A8").Value = 4
$ws.Range("B8").Value = "dragonball"
$ws.Range("C8").Value = "黑龙波"
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 40
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = 3

# --- 5. Grow the table to cover the two new rows. ------------------------
$lo.Resize($ws.Range("A3:G8"))

# --- 6. Match the author's final selection. -------------------------------
$ws.Range("G7").Select()
